# "This file replaced the User Reqmt Prioritization file"
#
# Updates the Smithgall Woods "System Requirements List" worksheet (Sheet1 /
# ActiveSheet):
#   - Row 5's "Web App Features" cell (B5) is renamed from the generic
#     "Website Pages (Site navigation & content presentation)" label to the
#     more specific "Website Pages (Content presentation)".
#   - Row 5's two example requirement cells (E5 / F5) get reworded/expanded
#     text (more explicit about contrast background & hover behaviour, and
#     about line-spacing / UX rationale).
#   - The frozen-pane / selection state on the sheet changes from a simple
#     row freeze to a combined row+column freeze anchored at D5, with the
#     final selection sitting on E5 in the bottom-right pane.
#   - Column widths for B, E and F are nudged slightly to better fit the
#     revised text, and row 5 grows taller to accommodate the longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Cell content updates (row 5) --------------------------------------
$ws.Range("B5").Value = "Website Pages (Content presentation)"
$ws.Range("E5").Value = "Each site page shall use relevant articles typed in darker font with higher contrast background & displayed next to their corresponding photos.  Currently, photos are placed on a page without corresponding label or reference except when hovering over photos.  Light font is used on a light gray background which can be difficult to read."
$ws.Range("F5").Value = "Each page shall use consistent line spacing, font type/size, and bold-type font that provides better user experience:  double-line spacing between paragraphs, after page titles, and after subheadings; 26pt bold font for page titles and 20pt bold font for subheadings.  Currently, some pages have inconsistent formatting."

# --- Row height grows to fit the longer wrapped text --------------------
$ws.Rows.Item(5).RowHeight = 189

# --- Column width tweaks ---------------------------------------------------
# (values chosen so the saved OOXML column width lands as close as possible
# to the target 32.5703125 / 31.28515625 / 30.5703125 given this host's
# pixel-snapped ColumnWidth setter)
$ws.Columns.Item(2).ColumnWidth = 31.5903
$ws.Columns.Item(5).ColumnWidth = 30.4252
$ws.Columns.Item(6).ColumnWidth = 29.5903

# --- Freeze panes: split at D5 (3 cols / 4 rows frozen), then leave the
#     live selection on E5 in the bottom-right pane -----------------------
[void]($excel.ActiveWindow.FreezePanes = $false)
$ws.Range("D5").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)
$ws.Range("E5").Select() | Out-Null
